$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.862.99"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.34"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.52"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3760"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07187"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8897"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08158"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.867.04"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.53"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.317"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.87"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.897.80"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.989"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.410"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.296"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.30"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.09"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.733"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.02"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.727"
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09148"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8153"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05007"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.177"
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.945"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6047"
$ws.Range("E36").Value = "  +5.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.212"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.616"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01952"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.071"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.638"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.934"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.22"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5104"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.952"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.634"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.71"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06072"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.27"
$ws.Range("E51").Value = "  -3.10%  "
